$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "9297354"
$ws.Range("D16").Value = "EDINSON JOSE RAMOS ZAMBRANO"
$ws.Range("E16").Value = "2404"
$ws.Range("F16").Value = 52000
$ws.Range("G16").Value = 1300000
$ws.Range("C17").Value = "1050953591"
$ws.Range("D17").Value = "YAIRA MARGARITA MUÑIZ GAMARRA"
$ws.Range("E17").Value = "2404"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000
$ws.Range("C18").Value = "1143385667"
$ws.Range("D18").Value = "DANILO HERNANDEZ RIVERA"
$ws.Range("E18").Value = "2404"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000
$ws.Range("C19").Value = "9297354"
$ws.Range("D19").Value = "EDINSON JOSE RAMOS ZAMBRANO"
$ws.Range("E19").Value = "2405"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000
$ws.Range("C20").Value = "1050953591"
$ws.Range("D20").Value = "YAIRA MARGARITA MUÑIZ GAMARRA"
$ws.Range("E20").Value = "2405"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1300000
$ws.Range("C21").Value = "1143385667"
$ws.Range("D21").Value = "DANILO HERNANDEZ RIVERA"
$ws.Range("E21").Value = "2405"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1300000
$ws.Range("C22").Value = "9297354"
$ws.Range("D22").Value = "EDINSON JOSE RAMOS ZAMBRANO"
$ws.Range("E22").Value = "2406"
$ws.Range("F22").Value = 52000
$ws.Range("G22").Value = 1300000
$ws.Range("C23").Value = "1050953591"
$ws.Range("D23").Value = "YAIRA MARGARITA MUÑIZ GAMARRA"
$ws.Range("E23").Value = "2406"
$ws.Range("F23").Value = 52000
$ws.Range("G23").Value = 1300000
$ws.Range("C24").Value = "1143385667"
$ws.Range("D24").Value = "DANILO HERNANDEZ RIVERA"
$ws.Range("E24").Value = "2406"
$ws.Range("F24").Value = 52000
$ws.Range("G24").Value = 1300000
$ws.Range("C25").Value = "9297354"
$ws.Range("D25").Value = "EDINSON JOSE RAMOS ZAMBRANO"
$ws.Range("E25").Value = "2407"
$ws.Range("F25").Value = 52000
$ws.Range("G25").Value = 1300000
$ws.Range("C26").Value = "1050953591"
$ws.Range("D26").Value = "YAIRA MARGARITA MUÑIZ GAMARRA"
$ws.Range("E26").Value = "2407"
$ws.Range("F26").Value = 52000
$ws.Range("G26").Value = 1300000
$ws.Range("C27").Value = "1143385667"
$ws.Range("D27").Value = "DANILO HERNANDEZ RIVERA"
$ws.Range("E27").Value = "2407"
$ws.Range("F27").Value = 52000
$ws.Range("G27").Value = 1300000
$ws.Range("C28").Value = "9297354"
$ws.Range("D28").Value = "EDINSON JOSE RAMOS ZAMBRANO"
$ws.Range("E28").Value = "2408"
$ws.Range("F28").Value = 52000
$ws.Range("G28").Value = 1300000
$ws.Range("C29").Value = "1050953591"
$ws.Range("D29").Value = "YAIRA MARGARITA MUÑIZ GAMARRA"
$ws.Range("E29").Value = "2408"
$ws.Range("F29").Value = 52000
$ws.Range("G29").Value = 1300000
$ws.Range("C30").Value = "1143385667"
$ws.Range("D30").Value = "DANILO HERNANDEZ RIVERA"
$ws.Range("E30").Value = "2408"
$ws.Range("F30").Value = 52000
$ws.Range("G30").Value = 1300000
$ws.Range("C31").Value = "9297354"
$ws.Range("D31").Value = "EDINSON JOSE RAMOS ZAMBRANO"
$ws.Range("E31").Value = "2409"
$ws.Range("F31").Value = 52000
$ws.Range("G31").Value = 1300000
$ws.Range("C32").Value = "1050953591"
$ws.Range("D32").Value = "YAIRA MARGARITA MUÑIZ GAMARRA"
$ws.Range("E32").Value = "2409"
$ws.Range("F32").Value = 52000
$ws.Range("G32").Value = 1300000
$ws.Range("C33").Value = "1143385667"
$ws.Range("D33").Value = "DANILO HERNANDEZ RIVERA"
$ws.Range("E33").Value = "2409"
$ws.Range("F33").Value = 52000
$ws.Range("G33").Value = 1300000
$ws.Range("C34").Value = "9297354"
$ws.Range("D34").Value = "EDINSON JOSE RAMOS ZAMBRANO"
$ws.Range("E34").Value = "2410"
$ws.Range("F34").Value = 15600
$ws.Range("G34").Value = 1300000
$ws.Range("C35").Value = "1050953591"
$ws.Range("D35").Value = "YAIRA MARGARITA MUÑIZ GAMARRA"
$ws.Range("E35").Value = "2410"
$ws.Range("F35").Value = 15600
$ws.Range("G35").Value = 1300000
$ws.Range("C36").Value = "1143385667"
$ws.Range("D36").Value = "DANILO HERNANDEZ RIVERA"
$ws.Range("E36").Value = "2410"
$ws.Range("F36").Value = 15600
$ws.Range("G36").Value = 1300000
